# This script applies a row-wise permutation of the "weekly" price data
# (columns D, J, K, L, M, P) across rows 2-22 of the active sheet, as
# described by the commit "Fruta / hortaliza, semanal".
#
# For each destination row, the values of D/J/K/L/M/P are replaced with
# the values that (before this edit) lived in a different ("source") row.
# Row 8 keeps its own original values (maps to itself).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (values to copy FROM, as they existed
# before any modification)
$map = @{
    2  = 4
    3  = 13
    4  = 19
    5  = 20
    6  = 18
    7  = 21
    8  = 8
    9  = 14
    10 = 17
    11 = 22
    12 = 9
    13 = 11
    14 = 2
    15 = 3
    16 = 6
    17 = 10
    18 = 12
    19 = 16
    20 = 5
    21 = 7
    22 = 15
}

$cols = @(4, 10, 11, 12, 13, 16)  # D, J, K, L, M, P

# Snapshot all original values first, since sources and destinations overlap.
# Value2 is used (rather than Value) because it yields plain CLR primitives
# (numbers/strings) that survive being stashed in a hashtable and written
# back later, whereas Value can return a COM Variant wrapper in this
# environment.
$original = @{}
foreach ($r in 2..22) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $original[$r] = $rowVals
}

# Now write the permuted values.
foreach ($r in 2..22) {
    $src = $map[$r]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value2 = $original[$src][$c]
    }
}
